$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 32 (shifts existing rows 32:85 down to 33:86)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with data (same shape as the rows around it)
$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(32, 3).Value = "Bíobío"
$ws.Cells.Item(32, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(32, 6).Value = 100112012
$ws.Cells.Item(32, 7).Value = "Espinaca"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 40
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 7500
$ws.Cells.Item(32, 13).Value = 7250
$ws.Cells.Item(32, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(32, 15).Value = "Región Metropolitana"
$ws.Cells.Item(32, 16).Value = 725
$ws.Cells.Item(32, 17).Value = 10
$ws.Cells.Item(32, 18).Value = "Hortaliza"

Write-Host "Done"
